$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D"  = 0.08869999999999999
    "G"  = 0.08419337316675719
    "H"  = 0.08419337316675719
    "I"  = 0.05111352525801195
    "J"  = 0.03297323492134496
    "K"  = 6.58
    "L"  = 0.03574144486692015
    "U"  = 12.1
    "V"  = 0.03705972434915773
    "W"  = 0.02811965811965812
    "X"  = 0.07693172573644817
    "Y"  = -0.04881206761679005
    "Z"  = 0.8410232983097304
    "AA" = 0.02773125878949112
    "AB" = 0.07693172573644817
    "AC" = -0.04920046694695705
    "AG" = -12.1
    "AJ" = -0.03848600508905853
    "AK" = -0.05309346204475647
    "AP" = -1.217303822937626
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
